$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (the "Förändrad" date column) for rows 2 through 16
# from 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244)
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = 45244
    }
}
